$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number but must remain text
# (matches the source data which stores prices/percentages as strings),
# so force a Text number format before assigning to avoid Excel
# auto-converting the literal into a numeric value.

$ws.Range("D2").Value = "69.699.22"

$ws.Range("D3").Value = "3.700.44"
$ws.Range("E3").Value = "  +0.61%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "671.85"
$ws.Range("E5").Value = "  -1.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.70"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.87%  "

$ws.Range("E9").Value = "  +0.49%  "

$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("E11").Value = "  +1.98%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000235"
$ws.Range("E12").Value = "  +1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.85"
$ws.Range("E13").Value = "  +1.86%  "

$ws.Range("D14").Value = "3.725.00"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("D15").Value = "69.689.28"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("E17").Value = "  +2.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.50"
$ws.Range("E18").Value = "  +1.82%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "474.85"
$ws.Range("E19").Value = "  +1.09%  "

$ws.Range("E21").Value = "  +0.82%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.48"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("D23").Value = "3.848.99"
$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("E24").Value = "  +3.70%  "

$ws.Range("E25").Value = "  -0.05%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.95"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.14"
$ws.Range("E27").Value = "  -0.55%  "

$ws.Range("E28").Value = "  -0.52%  "

$ws.Range("E29").Value = "  +0.99%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.02"
$ws.Range("E30").Value = "  +1.76%  "

$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  +4.31%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.89"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.996"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").Value = "3.690.17"
$ws.Range("E35").Value = "  +1.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.57"
$ws.Range("E36").Value = "  +5.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.24"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  +1.35%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "172.24"
$ws.Range("E42").Value = "  +3.66%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("E45").Value = "  +2.93%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000281"
$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.90"
$ws.Range("E47").Value = "  +2.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  -1.27%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("E49").Value = "  -1.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.90"
$ws.Range("E50").Value = "  +1.45%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.267"
$ws.Range("E51").Value = "  +0.83%  "
